# Auto-generated COM-interop script reproducing the workbook edit.
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Sheet "Pages produit": append product row 6
# ---------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Pages produit")
$ws1.Cells.Item(6, 1).Value = 'WESTERN, BOTTE SANTIAG VELOURS TAUPE'
$ws1.Cells.Item(6, 2).Value = 'Femme'
$ws1.Cells.Item(6, 3).Value = 'Shoes'
$ws1.Cells.Item(6, 4).Value = 'https://www.labottegardiane.com/products/western-paris-botte-santiag-velours-taupe#'
$ws1.Cells.Item(6, 5).Value = 7
# Row 4 has the exact same style pattern (s=9,10,10,9,10) needed for row 6
$ws1.Range("A4:E4").Copy()
$ws1.Range("A6:E6").PasteSpecial(-4122)
$ws1.Application.CutCopyMode = $false

# ---------------------------------------------------------------
# Sheet "Guides de taille": column A width 18 -> 20
# ---------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Guides de taille")
$ws2.Columns.Item(1).ColumnWidth = 19.17

# ---------------------------------------------------------------
# Pre-format the two numeric-looking-text blocks as Text so that
# values such as "39", "7.5" are stored as strings (matching the
# source, which keeps shoe sizes as inlineStr) instead of numbers.
# ---------------------------------------------------------------
$ws2.Range("C49:S53").NumberFormat = "@"
$ws2.Range("C59:S63").NumberFormat = "@"

# ---------------------------------------------------------------
# New "Guide de taille" block 6 (rows 46-54) and block 7 (rows 56-64)
# ---------------------------------------------------------------
# row 46
$ws2.Cells.Item(46, 1).Value = 'Guide de taille'
$ws2.Cells.Item(46, 2).Value = 6
$ws2.Cells.Item(46, 3).Value = 'URL'
$ws2.Cells.Item(46, 4).Value = 'https://www.labottegardiane.com/products/western-paris-botte-santiag-velours-taupe#'
# row 48
$ws2.Cells.Item(48, 1).Value = 'Systemes metriques'
$ws2.Cells.Item(48, 3).Value = 'Taille 1'
$ws2.Cells.Item(48, 4).Value = 'Taille 2'
$ws2.Cells.Item(48, 5).Value = 'Taille 3'
$ws2.Cells.Item(48, 6).Value = 'Taille 4'
$ws2.Cells.Item(48, 7).Value = 'Taille 5'
$ws2.Cells.Item(48, 8).Value = 'Taille 6'
$ws2.Cells.Item(48, 9).Value = 'Taille 7'
$ws2.Cells.Item(48, 10).Value = 'Taille 8'
$ws2.Cells.Item(48, 11).Value = 'Taille 9'
$ws2.Cells.Item(48, 12).Value = 'Taille 10'
$ws2.Cells.Item(48, 13).Value = 'Taille 11'
$ws2.Cells.Item(48, 14).Value = 'Taille 12'
$ws2.Cells.Item(48, 15).Value = 'Taille 13'
$ws2.Cells.Item(48, 16).Value = 'Taille 14'
$ws2.Cells.Item(48, 17).Value = 'Taille 15'
$ws2.Cells.Item(48, 18).Value = 'Taille 16'
$ws2.Cells.Item(48, 19).Value = 'Taille 17'
# row 49
$ws2.Cells.Item(49, 1).Value = 'La Bottega Gardiane'
$ws2.Cells.Item(49, 2).Value = 'La Bottega Gardiane'
$ws2.Cells.Item(49, 3).Value = '39'
$ws2.Cells.Item(49, 4).Value = '39.5'
$ws2.Cells.Item(49, 5).Value = '40'
$ws2.Cells.Item(49, 6).Value = '40.5'
$ws2.Cells.Item(49, 7).Value = '41'
$ws2.Cells.Item(49, 8).Value = '41.5'
$ws2.Cells.Item(49, 9).Value = '42'
$ws2.Cells.Item(49, 10).Value = '42.5'
$ws2.Cells.Item(49, 11).Value = '43'
$ws2.Cells.Item(49, 12).Value = '43.5'
$ws2.Cells.Item(49, 13).Value = '44'
$ws2.Cells.Item(49, 14).Value = '44.5'
$ws2.Cells.Item(49, 15).Value = '45'
$ws2.Cells.Item(49, 16).Value = '46'
$ws2.Cells.Item(49, 17).Value = '47'
$ws2.Cells.Item(49, 18).Value = '48'
$ws2.Cells.Item(49, 19).Value = '49'
# row 50
$ws2.Cells.Item(50, 1).Value = 'Europe'
$ws2.Cells.Item(50, 2).Value = 'EU'
$ws2.Cells.Item(50, 3).Value = '39'
$ws2.Cells.Item(50, 4).Value = '39.5'
$ws2.Cells.Item(50, 5).Value = '40'
$ws2.Cells.Item(50, 6).Value = '40.5'
$ws2.Cells.Item(50, 7).Value = '41'
$ws2.Cells.Item(50, 8).Value = '41.5'
$ws2.Cells.Item(50, 9).Value = '42'
$ws2.Cells.Item(50, 10).Value = '42.5'
$ws2.Cells.Item(50, 11).Value = '43'
$ws2.Cells.Item(50, 12).Value = '43.5'
$ws2.Cells.Item(50, 13).Value = '44'
$ws2.Cells.Item(50, 14).Value = '44.5'
$ws2.Cells.Item(50, 15).Value = '45'
$ws2.Cells.Item(50, 16).Value = '46'
$ws2.Cells.Item(50, 17).Value = '47'
$ws2.Cells.Item(50, 18).Value = '48'
$ws2.Cells.Item(50, 19).Value = '49'
# row 51
$ws2.Cells.Item(51, 1).Value = 'Royaume-Uni'
$ws2.Cells.Item(51, 2).Value = 'UK'
$ws2.Cells.Item(51, 3).Value = '5.5'
$ws2.Cells.Item(51, 4).Value = '6'
$ws2.Cells.Item(51, 5).Value = '6.5'
$ws2.Cells.Item(51, 6).Value = '7'
$ws2.Cells.Item(51, 7).Value = '7.5'
$ws2.Cells.Item(51, 8).Value = '7.5'
$ws2.Cells.Item(51, 9).Value = '8'
$ws2.Cells.Item(51, 10).Value = '8.5'
$ws2.Cells.Item(51, 11).Value = '9'
$ws2.Cells.Item(51, 12).Value = '9.5'
$ws2.Cells.Item(51, 13).Value = '9.5'
$ws2.Cells.Item(51, 14).Value = '10'
$ws2.Cells.Item(51, 15).Value = '10.5'
$ws2.Cells.Item(51, 16).Value = '11'
$ws2.Cells.Item(51, 17).Value = '12'
$ws2.Cells.Item(51, 18).Value = '13'
$ws2.Cells.Item(51, 19).Value = '14'
# row 52
$ws2.Cells.Item(52, 1).Value = 'Etats-Unis'
$ws2.Cells.Item(52, 2).Value = 'US'
$ws2.Cells.Item(52, 3).Value = '6'
$ws2.Cells.Item(52, 4).Value = '6.5'
$ws2.Cells.Item(52, 5).Value = '7'
$ws2.Cells.Item(52, 6).Value = '7.5'
$ws2.Cells.Item(52, 7).Value = '8'
$ws2.Cells.Item(52, 8).Value = '8'
$ws2.Cells.Item(52, 9).Value = '8.5'
$ws2.Cells.Item(52, 10).Value = '9'
$ws2.Cells.Item(52, 11).Value = '9.5'
$ws2.Cells.Item(52, 12).Value = '9.5'
$ws2.Cells.Item(52, 13).Value = '10'
$ws2.Cells.Item(52, 14).Value = '10.5'
$ws2.Cells.Item(52, 15).Value = '11'
$ws2.Cells.Item(52, 16).Value = '11.5'
$ws2.Cells.Item(52, 17).Value = '12.5'
$ws2.Cells.Item(52, 18).Value = '13.5'
$ws2.Cells.Item(52, 19).Value = '14.5'
# row 53
$ws2.Cells.Item(53, 1).Value = 'Italie'
$ws2.Cells.Item(53, 2).Value = 'IT'
$ws2.Cells.Item(53, 3).Value = '38'
$ws2.Cells.Item(53, 4).Value = '38.5'
$ws2.Cells.Item(53, 5).Value = '39'
$ws2.Cells.Item(53, 6).Value = '39.5'
$ws2.Cells.Item(53, 7).Value = '40'
$ws2.Cells.Item(53, 8).Value = '40.5'
$ws2.Cells.Item(53, 9).Value = '41'
$ws2.Cells.Item(53, 10).Value = '41.5'
$ws2.Cells.Item(53, 11).Value = '42'
$ws2.Cells.Item(53, 12).Value = '42.5'
$ws2.Cells.Item(53, 13).Value = '43'
$ws2.Cells.Item(53, 14).Value = '43.5'
$ws2.Cells.Item(53, 15).Value = '44'
$ws2.Cells.Item(53, 16).Value = '45'
$ws2.Cells.Item(53, 17).Value = '46'
$ws2.Cells.Item(53, 18).Value = '47'
$ws2.Cells.Item(53, 19).Value = '48'
# row 54
$ws2.Cells.Item(54, 1).Value = 'Longueur pied'
$ws2.Cells.Item(54, 3).Value = '25 cm'
$ws2.Cells.Item(54, 4).Value = '25.4 cm'
$ws2.Cells.Item(54, 5).Value = '25.7 cm'
$ws2.Cells.Item(54, 6).Value = '26 cm'
$ws2.Cells.Item(54, 7).Value = '26.4 cm'
$ws2.Cells.Item(54, 8).Value = '26.4 cm'
$ws2.Cells.Item(54, 9).Value = '26.7 cm'
$ws2.Cells.Item(54, 10).Value = '27.4 cm'
$ws2.Cells.Item(54, 11).Value = '27.7 cm'
$ws2.Cells.Item(54, 12).Value = '28 cm'
$ws2.Cells.Item(54, 13).Value = '28.4 cm'
$ws2.Cells.Item(54, 14).Value = '28.7 cm'
$ws2.Cells.Item(54, 15).Value = '29 cm'
$ws2.Cells.Item(54, 16).Value = '29.7 cm'
$ws2.Cells.Item(54, 17).Value = '30.4 cm'
$ws2.Cells.Item(54, 18).Value = '31 cm'
$ws2.Cells.Item(54, 19).Value = '31.7 cm'
# row 56
$ws2.Cells.Item(56, 1).Value = 'Guide de taille'
$ws2.Cells.Item(56, 2).Value = 7
$ws2.Cells.Item(56, 3).Value = 'URL'
$ws2.Cells.Item(56, 4).Value = 'https://www.labottegardiane.com/products/western-paris-botte-santiag-velours-taupe#'
# row 58
$ws2.Cells.Item(58, 1).Value = 'Systemes metriques'
$ws2.Cells.Item(58, 3).Value = 'Taille 1'
$ws2.Cells.Item(58, 4).Value = 'Taille 2'
$ws2.Cells.Item(58, 5).Value = 'Taille 3'
$ws2.Cells.Item(58, 6).Value = 'Taille 4'
$ws2.Cells.Item(58, 7).Value = 'Taille 5'
$ws2.Cells.Item(58, 8).Value = 'Taille 6'
$ws2.Cells.Item(58, 9).Value = 'Taille 7'
$ws2.Cells.Item(58, 10).Value = 'Taille 8'
$ws2.Cells.Item(58, 11).Value = 'Taille 9'
$ws2.Cells.Item(58, 12).Value = 'Taille 10'
$ws2.Cells.Item(58, 13).Value = 'Taille 11'
$ws2.Cells.Item(58, 14).Value = 'Taille 12'
$ws2.Cells.Item(58, 15).Value = 'Taille 13'
$ws2.Cells.Item(58, 16).Value = 'Taille 14'
$ws2.Cells.Item(58, 17).Value = 'Taille 15'
$ws2.Cells.Item(58, 18).Value = 'Taille 16'
$ws2.Cells.Item(58, 19).Value = 'Taille 17'
# row 59
$ws2.Cells.Item(59, 1).Value = 'La Bottega Gardiane'
$ws2.Cells.Item(59, 2).Value = 'La Bottega Gardiane'
$ws2.Cells.Item(59, 3).Value = '39'
$ws2.Cells.Item(59, 4).Value = '39.5'
$ws2.Cells.Item(59, 5).Value = '40'
$ws2.Cells.Item(59, 6).Value = '40.5'
$ws2.Cells.Item(59, 7).Value = '41'
$ws2.Cells.Item(59, 8).Value = '41.5'
$ws2.Cells.Item(59, 9).Value = '42'
$ws2.Cells.Item(59, 10).Value = '42.5'
$ws2.Cells.Item(59, 11).Value = '43'
$ws2.Cells.Item(59, 12).Value = '43.5'
$ws2.Cells.Item(59, 13).Value = '44'
$ws2.Cells.Item(59, 14).Value = '44.5'
$ws2.Cells.Item(59, 15).Value = '45'
$ws2.Cells.Item(59, 16).Value = '46'
$ws2.Cells.Item(59, 17).Value = '47'
$ws2.Cells.Item(59, 18).Value = '48'
$ws2.Cells.Item(59, 19).Value = '49'
# row 60
$ws2.Cells.Item(60, 1).Value = 'Europe'
$ws2.Cells.Item(60, 2).Value = 'EU'
$ws2.Cells.Item(60, 3).Value = '39'
$ws2.Cells.Item(60, 4).Value = '39.5'
$ws2.Cells.Item(60, 5).Value = '40'
$ws2.Cells.Item(60, 6).Value = '40.5'
$ws2.Cells.Item(60, 7).Value = '41'
$ws2.Cells.Item(60, 8).Value = '41.5'
$ws2.Cells.Item(60, 9).Value = '42'
$ws2.Cells.Item(60, 10).Value = '42.5'
$ws2.Cells.Item(60, 11).Value = '43'
$ws2.Cells.Item(60, 12).Value = '43.5'
$ws2.Cells.Item(60, 13).Value = '44'
$ws2.Cells.Item(60, 14).Value = '44.5'
$ws2.Cells.Item(60, 15).Value = '45'
$ws2.Cells.Item(60, 16).Value = '46'
$ws2.Cells.Item(60, 17).Value = '47'
$ws2.Cells.Item(60, 18).Value = '48'
$ws2.Cells.Item(60, 19).Value = '49'
# row 61
$ws2.Cells.Item(61, 1).Value = 'Royaume-Uni'
$ws2.Cells.Item(61, 2).Value = 'UK'
$ws2.Cells.Item(61, 3).Value = '5.5'
$ws2.Cells.Item(61, 4).Value = '6'
$ws2.Cells.Item(61, 5).Value = '6.5'
$ws2.Cells.Item(61, 6).Value = '7'
$ws2.Cells.Item(61, 7).Value = '7.5'
$ws2.Cells.Item(61, 8).Value = '7.5'
$ws2.Cells.Item(61, 9).Value = '8'
$ws2.Cells.Item(61, 10).Value = '8.5'
$ws2.Cells.Item(61, 11).Value = '9'
$ws2.Cells.Item(61, 12).Value = '9.5'
$ws2.Cells.Item(61, 13).Value = '9.5'
$ws2.Cells.Item(61, 14).Value = '10'
$ws2.Cells.Item(61, 15).Value = '10.5'
$ws2.Cells.Item(61, 16).Value = '11'
$ws2.Cells.Item(61, 17).Value = '12'
$ws2.Cells.Item(61, 18).Value = '13'
$ws2.Cells.Item(61, 19).Value = '14'
# row 62
$ws2.Cells.Item(62, 1).Value = 'Etats-Unis'
$ws2.Cells.Item(62, 2).Value = 'US'
$ws2.Cells.Item(62, 3).Value = '6'
$ws2.Cells.Item(62, 4).Value = '6.5'
$ws2.Cells.Item(62, 5).Value = '7'
$ws2.Cells.Item(62, 6).Value = '7.5'
$ws2.Cells.Item(62, 7).Value = '8'
$ws2.Cells.Item(62, 8).Value = '8'
$ws2.Cells.Item(62, 9).Value = '8.5'
$ws2.Cells.Item(62, 10).Value = '9'
$ws2.Cells.Item(62, 11).Value = '9.5'
$ws2.Cells.Item(62, 12).Value = '9.5'
$ws2.Cells.Item(62, 13).Value = '10'
$ws2.Cells.Item(62, 14).Value = '10.5'
$ws2.Cells.Item(62, 15).Value = '11'
$ws2.Cells.Item(62, 16).Value = '11.5'
$ws2.Cells.Item(62, 17).Value = '12.5'
$ws2.Cells.Item(62, 18).Value = '13.5'
$ws2.Cells.Item(62, 19).Value = '14.5'
# row 63
$ws2.Cells.Item(63, 1).Value = 'Italie'
$ws2.Cells.Item(63, 2).Value = 'IT'
$ws2.Cells.Item(63, 3).Value = '38'
$ws2.Cells.Item(63, 4).Value = '38.5'
$ws2.Cells.Item(63, 5).Value = '39'
$ws2.Cells.Item(63, 6).Value = '39.5'
$ws2.Cells.Item(63, 7).Value = '40'
$ws2.Cells.Item(63, 8).Value = '40.5'
$ws2.Cells.Item(63, 9).Value = '41'
$ws2.Cells.Item(63, 10).Value = '41.5'
$ws2.Cells.Item(63, 11).Value = '42'
$ws2.Cells.Item(63, 12).Value = '42.5'
$ws2.Cells.Item(63, 13).Value = '43'
$ws2.Cells.Item(63, 14).Value = '43.5'
$ws2.Cells.Item(63, 15).Value = '44'
$ws2.Cells.Item(63, 16).Value = '45'
$ws2.Cells.Item(63, 17).Value = '46'
$ws2.Cells.Item(63, 18).Value = '47'
$ws2.Cells.Item(63, 19).Value = '48'
# row 64
$ws2.Cells.Item(64, 1).Value = 'Longueur pied'
$ws2.Cells.Item(64, 3).Value = '25 cm'
$ws2.Cells.Item(64, 4).Value = '25.4 cm'
$ws2.Cells.Item(64, 5).Value = '25.7 cm'
$ws2.Cells.Item(64, 6).Value = '26 cm'
$ws2.Cells.Item(64, 7).Value = '26.4 cm'
$ws2.Cells.Item(64, 8).Value = '26.4 cm'
$ws2.Cells.Item(64, 9).Value = '26.7 cm'
$ws2.Cells.Item(64, 10).Value = '27.4 cm'
$ws2.Cells.Item(64, 11).Value = '27.7 cm'
$ws2.Cells.Item(64, 12).Value = '28 cm'
$ws2.Cells.Item(64, 13).Value = '28.4 cm'
$ws2.Cells.Item(64, 14).Value = '28.7 cm'
$ws2.Cells.Item(64, 15).Value = '29 cm'
$ws2.Cells.Item(64, 16).Value = '29.7 cm'
$ws2.Cells.Item(64, 17).Value = '30.4 cm'
$ws2.Cells.Item(64, 18).Value = '31 cm'
$ws2.Cells.Item(64, 19).Value = '31.7 cm'

# ---------------------------------------------------------------
# Apply formatting (fill/font/alignment) by copying style from the
# matching template cells in the existing "Guide de taille" block 5
# (rows 37-44), which has the identical style pattern.
# ---------------------------------------------------------------
# row 46 formatting
$ws2.Range("A1").Copy()
$ws2.Range("A46").PasteSpecial(-4122)
$ws2.Range("B1").Copy()
$ws2.Range("B46").PasteSpecial(-4122)
$ws2.Range("C1").Copy()
$ws2.Range("C46").PasteSpecial(-4122)
$ws2.Range("D19").Copy()
$ws2.Range("D46").PasteSpecial(-4122)
# row 48 formatting
$ws2.Range("A1").Copy()
$ws2.Range("A48").PasteSpecial(-4122)
$ws2.Range("B3").Copy()
$ws2.Range("B48").PasteSpecial(-4122)
$ws2.Range("A1").Copy()
$ws2.Range("C48:S48").PasteSpecial(-4122)
# row 49 formatting
$ws2.Range("A4").Copy()
$ws2.Range("A49").PasteSpecial(-4122)
$ws2.Range("B4").Copy()
$ws2.Range("B49").PasteSpecial(-4122)
$ws2.Range("A4").Copy()
$ws2.Range("C49:S49").PasteSpecial(-4122)
# row 50 formatting
$ws2.Range("A4").Copy()
$ws2.Range("A50").PasteSpecial(-4122)
$ws2.Range("B4").Copy()
$ws2.Range("B50").PasteSpecial(-4122)
$ws2.Range("A4").Copy()
$ws2.Range("C50:S50").PasteSpecial(-4122)
# row 51 formatting
$ws2.Range("A4").Copy()
$ws2.Range("A51").PasteSpecial(-4122)
$ws2.Range("B4").Copy()
$ws2.Range("B51").PasteSpecial(-4122)
$ws2.Range("A4").Copy()
$ws2.Range("C51:S51").PasteSpecial(-4122)
# row 52 formatting
$ws2.Range("A4").Copy()
$ws2.Range("A52").PasteSpecial(-4122)
$ws2.Range("B4").Copy()
$ws2.Range("B52").PasteSpecial(-4122)
$ws2.Range("A4").Copy()
$ws2.Range("C52:S52").PasteSpecial(-4122)
# row 53 formatting
$ws2.Range("A4").Copy()
$ws2.Range("A53").PasteSpecial(-4122)
$ws2.Range("B4").Copy()
$ws2.Range("B53").PasteSpecial(-4122)
$ws2.Range("A4").Copy()
$ws2.Range("C53:S53").PasteSpecial(-4122)
# row 54 formatting
$ws2.Range("A4").Copy()
$ws2.Range("A54").PasteSpecial(-4122)
$ws2.Range("A4").Copy()
$ws2.Range("C54:S54").PasteSpecial(-4122)
# row 56 formatting
$ws2.Range("A1").Copy()
$ws2.Range("A56").PasteSpecial(-4122)
$ws2.Range("B1").Copy()
$ws2.Range("B56").PasteSpecial(-4122)
$ws2.Range("C1").Copy()
$ws2.Range("C56").PasteSpecial(-4122)
$ws2.Range("D19").Copy()
$ws2.Range("D56").PasteSpecial(-4122)
# row 58 formatting
$ws2.Range("A1").Copy()
$ws2.Range("A58").PasteSpecial(-4122)
$ws2.Range("B3").Copy()
$ws2.Range("B58").PasteSpecial(-4122)
$ws2.Range("A1").Copy()
$ws2.Range("C58:S58").PasteSpecial(-4122)
# row 59 formatting
$ws2.Range("A4").Copy()
$ws2.Range("A59").PasteSpecial(-4122)
$ws2.Range("B4").Copy()
$ws2.Range("B59").PasteSpecial(-4122)
$ws2.Range("A4").Copy()
$ws2.Range("C59:S59").PasteSpecial(-4122)
# row 60 formatting
$ws2.Range("A4").Copy()
$ws2.Range("A60").PasteSpecial(-4122)
$ws2.Range("B4").Copy()
$ws2.Range("B60").PasteSpecial(-4122)
$ws2.Range("A4").Copy()
$ws2.Range("C60:S60").PasteSpecial(-4122)
# row 61 formatting
$ws2.Range("A4").Copy()
$ws2.Range("A61").PasteSpecial(-4122)
$ws2.Range("B4").Copy()
$ws2.Range("B61").PasteSpecial(-4122)
$ws2.Range("A4").Copy()
$ws2.Range("C61:S61").PasteSpecial(-4122)
# row 62 formatting
$ws2.Range("A4").Copy()
$ws2.Range("A62").PasteSpecial(-4122)
$ws2.Range("B4").Copy()
$ws2.Range("B62").PasteSpecial(-4122)
$ws2.Range("A4").Copy()
$ws2.Range("C62:S62").PasteSpecial(-4122)
# row 63 formatting
$ws2.Range("A4").Copy()
$ws2.Range("A63").PasteSpecial(-4122)
$ws2.Range("B4").Copy()
$ws2.Range("B63").PasteSpecial(-4122)
$ws2.Range("A4").Copy()
$ws2.Range("C63:S63").PasteSpecial(-4122)
# row 64 formatting
$ws2.Range("A4").Copy()
$ws2.Range("A64").PasteSpecial(-4122)
$ws2.Range("A4").Copy()
$ws2.Range("C64:S64").PasteSpecial(-4122)
$ws2.Application.CutCopyMode = $false
